# Applies the "Atualizado por script em 19-12-2023 18:56" update to the
# costa-rica_primera-division_2023-2024 sheet:
#   1) Swap the match data (columns F:V) between existing rows 91 and 92.
#   2) Append 6 new match rows (134-139) with the same formatting as row 133.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap F91:V91 <-> F92:V92 -------------------------------------------
$rng91 = $ws.Range("F91:V91")
$rng92 = $ws.Range("F92:V92")
$vals91 = $rng91.Value2
$vals92 = $rng92.Value2
$rng91.Value2 = $vals92
$rng92.Value2 = $vals91

# --- 2) Append new rows 134-139 ---------------------------------------------
# Copy the formatting (styles) of row 133 down into the new rows first, so the
# bold/bordered index column (A) and the date column (E) keep their look.
$fmtSrc = $ws.Range("A133:V133")
$fmtSrc.Copy()
for ($r = 134; $r -le 139; $r++) {
    $dst = $ws.Range("A$r" + ":V$r")
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

$newRows = @(
    @(133, 'costa-rica', 'primera-division', '2023-2024', 45263.125, 'Herediano', 3, 'Alajuelense', 0, 1.99, '28/11/2023 11:12', 2.34, '03/12/2023 02:55', 3.52, '28/11/2023 11:12', 3.14, '03/12/2023 02:55', 3.72, '28/11/2023 11:12', 3.35, '03/12/2023 02:55', 'https://www.betexplorer.com/football/costa-rica/primera-division/herediano-alajuelense/0lhpQQ3k/'),
    @(134, 'costa-rica', 'primera-division', '2023-2024', 45263.75, 'Cartagines', 0, 'Saprissa', 2, 3.25, '28/11/2023 11:12', 3.46, '03/12/2023 17:58', 3.47, '28/11/2023 11:12', 3.39, '03/12/2023 17:58', 2.18, '28/11/2023 11:12', 2.17, '03/12/2023 17:58', 'https://www.betexplorer.com/football/costa-rica/primera-division/cartagines-saprissa/dGogO4Y1/'),
    @(135, 'costa-rica', 'primera-division', '2023-2024', 45270.125, 'Saprissa', 4, 'Cartagines', 0, 1.4, '06/12/2023 19:42', 1.55, '10/12/2023 02:30', 4.57, '06/12/2023 19:42', 4.13, '10/12/2023 02:59', 6.56, '06/12/2023 19:42', 6, '10/12/2023 02:59', 'https://www.betexplorer.com/football/costa-rica/primera-division/saprissa-cartagines/AJZzmlB2/'),
    @(136, 'costa-rica', 'primera-division', '2023-2024', 45271, 'Alajuelense', 1, 'Herediano', 0, 2.13, '06/12/2023 21:12', 1.99, '10/12/2023 23:59', 3.37, '06/12/2023 21:12', 3.43, '10/12/2023 23:59', 3.23, '06/12/2023 21:12', 3.98, '10/12/2023 23:59', 'https://www.betexplorer.com/football/costa-rica/primera-division/alajuelense-herediano/rDskPpJe/'),
    @(137, 'costa-rica', 'primera-division', '2023-2024', 45275.125, 'Herediano', 1, 'Saprissa', 2, 2.37, '13/12/2023 17:12', 2.26, '15/12/2023 02:56', 3.22, '13/12/2023 17:12', 3.15, '15/12/2023 02:55', 2.93, '13/12/2023 17:12', 3.66, '15/12/2023 02:38', 'https://www.betexplorer.com/football/costa-rica/primera-division/herediano-saprissa/CxsL5zMR/'),
    @(138, 'costa-rica', 'primera-division', '2023-2024', 45278, 'Saprissa', 1, 'Herediano', 0, 1.76, '15/12/2023 13:43', 2.1, '17/12/2023 23:58', 3.59, '15/12/2023 13:43', 3.34, '17/12/2023 23:58', 4.33, '15/12/2023 13:43', 3.69, '17/12/2023 23:58', 'https://www.betexplorer.com/football/costa-rica/primera-division/saprissa-herediano/nVXn1fUr/')
)

$startRow = 134
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value2 = $rowVals[$c]
    }
}
